$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.827.72"
$ws.Range("E2").Value = "  +0.20%  "

# Row 3
$ws.Range("D3").Value = "2.536.14"
$ws.Range("E3").Value = "  -0.38%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.70"
$ws.Range("E5").Value = "  +1.85%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.87"
$ws.Range("E6").Value = "  +4.00%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.577"
$ws.Range("E7").Value = "  +0.69%  "

# Row 8
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.546"
$ws.Range("E9").Value = "  -0.14%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.42"
$ws.Range("E10").Value = "  +0.96%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0823"
$ws.Range("E11").Value = "  +2.35%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.114"
$ws.Range("E12").Value = "  -0.41%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.59"
$ws.Range("E13").Value = "  -1.47%  "

# Row 14
$ws.Range("D14").Value = "2.926.33"

# Row 15
$ws.Range("D15").Value = "2.520.19"
$ws.Range("E15").Value = "  -1.77%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.07"
$ws.Range("E16").Value = "  +6.09%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.870"
$ws.Range("E17").Value = "  -0.83%  "

# Row 18
$ws.Range("D18").Value = "42.819.71"
$ws.Range("E18").Value = "  +0.17%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.20"
$ws.Range("E19").Value = "  +3.93%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0989"
$ws.Range("E20").Value = "  +0.94%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.55"
$ws.Range("E21").Value = "  +0.05%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.58"
$ws.Range("E22").Value = "  -0.14%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.40"
$ws.Range("E23").Value = "  +0.17%  "

# Row 24
$ws.Range("E24").Value = "  +0.27%  "

# Row 25
$ws.Range("E25").Value = "  -2.99%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.72"
$ws.Range("E26").Value = "  -4.23%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.17%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.29"
$ws.Range("E28").Value = "  +8.42%  "

# Row 29
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.19"
$ws.Range("E29").Value = "  +0.52%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.17"
$ws.Range("E30").Value = "  +3.41%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.20"
$ws.Range("E31").Value = "  +2.61%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.20"
$ws.Range("E32").Value = "  +3.29%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.39"
$ws.Range("E33").Value = "  +13.18%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.13"
$ws.Range("E34").Value = "  -0.98%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0796"
$ws.Range("E35").Value = "  +0.42%  "

# Row 36
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.30"
$ws.Range("E36").Value = "  -1.89%  "

# Row 37
$ws.Range("E37").Value = "  -4.66%  "

# Row 38
$ws.Range("E38").Value = "  +0.69%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.90"
$ws.Range("E39").Value = "  +7.04%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.120"
$ws.Range("E40").Value = "  +0.47%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.17"
$ws.Range("E41").Value = "  +31.29%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.42"
$ws.Range("E42").Value = "  +0.62%  "

# Row 43
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.101.02"
$ws.Range("E43").Value = "  +0.97%  "

# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.85"
$ws.Range("E44").Value = "  -0.30%  "

# Row 45
$ws.Range("E45").Value = "  -1.76%  "

# Row 46
$ws.Range("E46").Value = "  +0.03%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.53"
$ws.Range("E47").Value = "  +2.63%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.96"
$ws.Range("E48").Value = "  -0.91%  "

# Row 49
$ws.Range("D49").Value = "2.784.52"
$ws.Range("E49").Value = "  -0.23%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.57"
$ws.Range("E50").Value = "  +6.64%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.192"
$ws.Range("E51").Value = "  +1.86%  "
